$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue $ws "D2" "58.120.11"
Set-TextValue $ws "E2" "  -0.64%  "
Set-TextValue $ws "D3" "2.284.68"
Set-TextValue $ws "E3" "  +0.34%  "
Set-TextValue $ws "E4" "  -0.06%  "
Set-TextValue $ws "D5" "534.94"
Set-TextValue $ws "E5" "  -1.87%  "
Set-TextValue $ws "D6" "131.27"
Set-TextValue $ws "E6" "  +0.41%  "
Set-TextValue $ws "D7" "0.999"
Set-TextValue $ws "E7" "  -0.06%  "
Set-TextValue $ws "E8" "  +3.38%  "
Set-TextValue $ws "D9" "2.280.15"
Set-TextValue $ws "E9" "  +0.21%  "
Set-TextValue $ws "E10" "  -1.51%  "
Set-TextValue $ws "E11" "  -0.89%  "
Set-TextValue $ws "E12" "  +0.77%  "
Set-TextValue $ws "E13" "  -0.67%  "
Set-TextValue $ws "E14" "  -0.71%  "
Set-TextValue $ws "D15" "2.690.67"
Set-TextValue $ws "E15" "  +0.26%  "
Set-TextValue $ws "D16" "58.009.34"
Set-TextValue $ws "E16" "  -0.80%  "
Set-TextValue $ws "E17" "  -0.56%  "
Set-TextValue $ws "D18" "2.326.63"
Set-TextValue $ws "E18" "  +2.20%  "
Set-TextValue $ws "D19" "10.53"
Set-TextValue $ws "E19" "  -0.84%  "
Set-TextValue $ws "E20" "  -2.48%  "
Set-TextValue $ws "D21" "313.07"
Set-TextValue $ws "E21" "  -0.15%  "
Set-TextValue $ws "D22" "6.47"
Set-TextValue $ws "E22" "  +0.74%  "
Set-TextValue $ws "E23" "  -0.06%  "
Set-TextValue $ws "D24" "63.14"
Set-TextValue $ws "E25" "  -1.45%  "
Set-TextValue $ws "D26" "1.00"
Set-TextValue $ws "E26" "  +0.09%  "
Set-TextValue $ws "E27" "  -1.68%  "
Set-TextValue $ws "D28" "1.28"
Set-TextValue $ws "E28" "  -2.19%  "
Set-TextValue $ws "D29" "170.52"
Set-TextValue $ws "E29" "  -0.06%  "
Set-TextValue $ws "D30" "1.69"
Set-TextValue $ws "E30" "  -2.93%  "
Set-TextValue $ws "D31" "0.0₃0722"
Set-TextValue $ws "E31" "  +0.04%  "
Set-TextValue $ws "E32" "  +0.11%  "
Set-TextValue $ws "E33" "  -1.91%  "
Set-TextValue $ws "D34" "0.380"
Set-TextValue $ws "E34" "  -0.56%  "
Set-TextValue $ws "E35" "  -0.01%  "
Set-TextValue $ws "D36" "17.83"
Set-TextValue $ws "E36" "  +0.35%  "
Set-TextValue $ws "E37" "  -0.09%  "
Set-TextValue $ws "D38" "1.24"
Set-TextValue $ws "E39" "  -1.01%  "
Set-TextValue $ws "E40" "  -1.55%  "
Set-TextValue $ws "D41" "288.46"
Set-TextValue $ws "E41" "  -4.78%  "
Set-TextValue $ws "D42" "139.94"
Set-TextValue $ws "E42" "  -0.50%  "
Set-TextValue $ws "E43" "  -0.46%  "
Set-TextValue $ws "E44" "  +0.63%  "
Set-TextValue $ws "D45" "0.0495"
Set-TextValue $ws "E45" "  -0.44%  "
Set-TextValue $ws "D46" "0.553"
Set-TextValue $ws "E46" "  +0.45%  "
Set-TextValue $ws "D47" "18.08"
Set-TextValue $ws "E47" "  -1.31%  "
Set-TextValue $ws "E48" "  -1.78%  "
Set-TextValue $ws "D49" "10.95"
Set-TextValue $ws "E49" "  -0.66%  "
Set-TextValue $ws "E51" "  +1.23%  "
